$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet grows from 44 data rows (A1:D44) to 45 data rows (A1:D45): the
# "Z03_B04_P01" entry moves from the bottom (old row 44) up near the top (new
# row 8), a brand-new "Z04_B01_P02" entry is inserted after "Z04_B01_P01", and
# every row in between shifts down by one. Simplest + most robust: copy the
# format of the last existing data row down into the new row 45, then rewrite
# every data cell (A2:D45) with its final target value.
$ws.Range("A44:D44").Copy() | Out-Null
$ws.Range("A45:D45").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Final target data for rows 2-45 (columns A:D). Row 1 (header) is unchanged.
$data = @(
    ,@("Z01_B01_P01", "Z01_B01", "Armut begrenzen", "Limiting poverty")
    ,@("Z02_B01_P01", "Z02_B01", "In unseren Kulturlandschaften umweltverträglich produzieren", "Environmentally sound production in our cultivated landscapes")
    ,@("Z02_B02_P01", "Z02_B02", "Das Recht auf Nahrung weltweit verwirklichen", "Realising globally the right to food")
    ,@("Z03_B01_P01", "Z03_B01", "Länger gesund leben", "Living healthy longer")
    ,@("Z03_B02_P01", "Z03_B02", "Gesunde Umwelt erhalten", "Keeping the environment healthy")
    ,@("Z03_B03_P01", "Z03_B03", "Globale Gesundheitsarchitektur stärken", "Strengthening the global health architecture")
    ,@("Z03_B04_P01", "Z03_B04", "Unterschiede in der mittleren Lebenserwartung zwischen den Lebensräumen mit hoher bzw. niedriger sozioökonomischer Deprivation reduzieren", "XXX")
    ,@("Z04_B01_P01", "Z04_B01", "Bildung und Qualifikation kontinuierlich verbessern", "Continuously improving education and vocational training")
    ,@("Z04_B01_P02", "Z04_B01", "Verankerung von Bildung für nachhaltige Entwicklung im Bildungssystem vorantreiben", "XXX")
    ,@("Z04_B02_P01", "Z04_B02", "Vereinbarkeit von Familie und Beruf verbessern", "Improving the compatibility of work and family life")
    ,@("Z05_B01_P01", "Z05_B01", "Gleichstellung und partnerschaftliche Aufgabenteilung fördern", "Promoting equal opportunities in society")
    ,@("Z05_B01_P02", "Z05_B01", "Wirtschaftliche Teilhabe von Frauen global stärken", "Strengthening the economic participation of women globally")
    ,@("Z06_B01_P01", "Z06_B01", "Minderung der stofflichen Belastung von Gewässern", "Reducing the pollution of water with substances")
    ,@("Z06_B02_P01", "Z06_B02", "Besserer Zugang zu Trinkwasser- und Sanitärversorgung weltweit, höhere (sichere) Qualität", "Better access to drinking water and sanitation worldwide, higher (safer) quality")
    ,@("Z07_B01_P01", "Z07_B01", "Ressourcen sparsam und effizient nutzen", "Using resources economically and efficiently")
    ,@("Z07_B02_P01", "Z07_B02", "Zukunftsfähige Energieversorgung ausbauen", "Strengthening a sustainable energy supply")
    ,@("Z08_B01_P01", "Z08_B01", "Ressourcen sparsam und effizient nutzen", "Using resources economically and efficiently")
    ,@("Z08_B02_P01", "Z08_B02", "Staatsfinanzen konsolidieren – Generationengerechtigkeit schaffen", "Consolidating public finances – Creating intergenerational equity")
    ,@("Z08_B03_P01", "Z08_B03", "Gute Investitionsbedingungen schaffen – Wohlstand dauerhaft erhalten", "Creating favourable investment conditions – Securing long-term prosperity")
    ,@("Z08_B04_P01", "Z08_B04", "Wirtschaftsleistung umwelt- und sozialverträglich steigern", "Combining greater economic output with environmental and social responsibility")
    ,@("Z08_B05_P01", "Z08_B05", "Beschäftigungsniveau steigern", "Boosting employment levels")
    ,@("Z08_B06_P01", "Z08_B06", "Menschenwürdige Arbeit weltweit ermöglichen", "Enabling decent work worldwide")
    ,@("Z09_B01_P01", "Z09_B01", "Zukunft mit neuen Lösungen nachhaltig gestalten", "Shaping the future with new solutions")
    ,@("Z10_B01_P01", "Z10_B01", "Schulische Bildungserfolge von Ausländern in Deutschland verbessern", "Improving educational success of foreigners in German schools")
    ,@("Z10_B02_P01", "Z10_B02", "Zu große Ungleichheit innerhalb Deutschlands verhindern", "Preventing excessive inequality within Germany")
    ,@("Z11_B01_P01", "Z11_B01", "Flächen nachhaltig nutzen", "Using land sustainably")
    ,@("Z11_B02_P01", "Z11_B02", "Mobilität sichern – Umwelt schonen", "Guaranteeing mobility – Protecting the environment")
    ,@("Z11_B03_P01", "Z11_B03", "Bezahlbarer Wohnraum für alle", "Affordable housing for all")
    ,@("Z11_B04_P01", "Z11_B04", "Zugang zum Kulturerbe verbessern", "Improving access to cultural heritage")
    ,@("Z12_B01_P01", "Z12_B01", "Konsum umwelt- und sozialverträglich gestalten", "Making consumption environmentally and socially compatible")
    ,@("Z12_B02_P01", "Z12_B02", "Anteil nachhaltiger Produktion stetig erhöhen", "Increasing the proportion of sustainable production continuously")
    ,@("Z12_B03_P01", "Z12_B03", "Vorbildwirkung der öffentlichen Hand für nachhaltige öffentliche Beschaffung verwirklichen", "Giving shape to the public sector’s exemplary role in sustainable procurement")
    ,@("Z13_B01_P01", "Z13_B01", "Treibhausgase reduzieren", "Reducing greenhouse gases")
    ,@("Z13_B01_P02", "Z13_B01", "Beitrag zur internationalen Klimafinanzierung leisten", "Germany's contribution to international climate finance")
    ,@("Z14_B01_P01", "Z14_B01", "Meere und Meeresressourcen schützen und nachhaltig nutzen", "Protecting and sustainably using oceans and marine resources")
    ,@("Z15_B01_P01", "Z15_B01", "Arten erhalten – Lebensräume schützen", "Conserving species – Protecting habitats")
    ,@("Z15_B02_P01", "Z15_B02", "Ökosysteme schützen, Ökosystemleistungen erhalten und Lebensräume bewahren", "Protecting ecosystems, conserving ecosystem services and preserving habitats")
    ,@("Z15_B02_P02", "Z15_B02", "Weltweit Entwaldung vermeiden und Böden schützen", "Preventing deforestation and protecting soils world-wide")
    ,@("Z16_B01_P01", "Z16_B01", "Persönliche Sicherheit weiter erhöhen", "Further increasing personal security")
    ,@("Z16_B02_P01", "Z16_B02", "Praktische Maßnahmen zur Bekämpfung der Proliferation, insbesondere von Kleinwaffen, ergreifen", "Taking practical action to combat proliferation, especially of small arms")
    ,@("Z16_B03_P01", "Z16_B03", "Korruption bekämpfen", "Combating corruption")
    ,@("Z17_B01_P01", "Z17_B01", "Nachhaltige Entwicklung unterstützen", "Supporting sustainable development")
    ,@("Z17_B02_P01", "Z17_B02", "Wissen international vermitteln", "Sharing knowledge internationally")
    ,@("Z17_B03_P01", "Z17_B03", "Handelschancen der Entwicklungsländer verbessern", "Improving trade opportunities for developing countries")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowVals[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowVals[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowVals[2]
    $ws.Cells.Item($rowNum, 4).Value = $rowVals[3]
}